# Auto-generated data-driven edit script for COMS worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert new columns D (Corequisites), E (Concurrent), F (Recommended);
# shift old 'Terms Typically Offered' header from D1 to G1
$ws.Cells.Item(1, 4).Value = 'Corequisites'
$ws.Cells.Item(1, 5).Value = 'Concurrent'
$ws.Cells.Item(1, 6).Value = 'Recommended'
$ws.Cells.Item(1, 7).Value = 'Terms Typically Offered'

# Per-row data: for rows 2-47, the old 'Terms Typically Offered' value (col D)
# moves to col G. New col D (Corequisites) and E (Concurrent) are 'NA' for every row.
# Col F (Recommended) is 'NA' unless the row's old Prerequisites text (col C) embedded
# a 'Recommended: ...' clause, which is extracted into F and stripped from C.
$rows = @(
    @{ Row = 2; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F,W,SP,SU' }
    @{ Row = 3; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 4; NewC = 'Completion of GE Area A1 with a grade of C- or better.'; D = 'NA'; E = 'NA'; F = 'Completion of GE Area A2.'; G = 'F, W, SP ' }
    @{ Row = 5; NewC = 'Completion of GE Area A1 with a grade of C- or better.'; D = 'NA'; E = 'NA'; F = 'Completion of GE Area A2.'; G = 'W, SP ' }
    @{ Row = 6; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 7; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 8; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 9; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 10; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 11; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F,W,SP,SU' }
    @{ Row = 12; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 13; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 14; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 15; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 16; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W, SP' }
    @{ Row = 17; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 18; NewC = 'COMS 311 and STAT 217, junior standing; for COMS majors only.'; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 19; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 20; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F,W,SP,SU' }
    @{ Row = 21; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 22; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 23; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 24; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 25; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 26; NewC = 'Junior standing, COMS 330.'; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 27; NewC = 'COMS 250.'; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 28; NewC = $null; D = 'NA'; E = 'NA'; F = 'COMS 218 or JOUR 218.'; G = 'F, W, Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D. ' }
    @{ Row = 29; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W' }
    @{ Row = 30; NewC = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'; D = 'NA'; E = 'NA'; F = 'COMS 218 or POLS 112.'; G = 'F, SU ' }
    @{ Row = 31; NewC = 'Junior standing and completion of GE Area A with grades of C- or better.'; D = 'NA'; E = 'NA'; F = 'Completion of GE Area B2 or B3.'; G = 'W ' }
    @{ Row = 32; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 33; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 34; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'Junior standing and completion of GE Area A with grades of C- or better.' }
    @{ Row = 35; NewC = 'Junior standing.'; D = 'NA'; E = 'NA'; F = 'COMS 213 or COMS 301.'; G = 'W ' }
    @{ Row = 36; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 37; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 38; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 39; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 40; NewC = 'Completion of GE Area A with grades of C- or better; and COMS/PSY 212.'; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 41; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W' }
    @{ Row = 42; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F,W,SP,SU' }
    @{ Row = 43; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W' }
    @{ Row = 44; NewC = 'Communication Studies major; COMS 312; COMS 332; and junior standing.'; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 45; NewC = 'COMS 460; for COMS majors only.'; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 46; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 47; NewC = $null; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F,W,SP,SU' }
)

foreach ($r in $rows) {
    if ($null -ne $r.NewC) {
        $ws.Cells.Item($r.Row, 3).Value = $r.NewC
    }
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

